# Update countries & provincias Spain
#
# The source "Pais" sheet lists one country per row (rows 4-216), sorted
# descending by "Casos totales" (column B). A handful of countries got
# refreshed case counts; refreshing them changes their sort rank, so the
# whole table is re-sorted afterwards. The "last updated" timestamp in A1
# is bumped too.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Bump the "last updated" timestamp (row 1) ------------------------
$ws.Range("A1").Value = "Datos actualizados a 26 de Abril de 2020 a las 04:22"

# --- 2. Refresh the handful of countries with new case numbers -----------
# Columns: B=Casos totales, C=Nuevos casos, D=Casos activos, E=Recuperados,
#          F=Casos criticos, G=Muertes hoy, H=Muertes

function Set-CountryRow($rowNum, $total, $nuevos, $activos, $recuperados, $criticos, $muertesHoy, $muertes) {
    $ws.Cells.Item($rowNum, 2).Value = $total
    $ws.Cells.Item($rowNum, 3).Value = $nuevos
    $ws.Cells.Item($rowNum, 4).Value = $activos
    $ws.Cells.Item($rowNum, 5).Value = $recuperados
    $ws.Cells.Item($rowNum, 6).Value = $criticos
    $ws.Cells.Item($rowNum, 7).Value = $muertesHoy
    $ws.Cells.Item($rowNum, 8).Value = $muertes
}

# Bolivia (row 92 before the re-sort)
Set-CountryRow 92 866 59 54 766 3 2 46

# Estado de Palestina (row 106 before the re-sort)
Set-CountryRow 106 342 0 92 248 0 0 2

# Sri Lanka (row 108 before the re-sort)
Set-CountryRow 108 460 8 118 335 2 0 7

# Birmania (row 137 before the re-sort) - no reorder needed, still updates
Set-CountryRow 137 146 2 10 131 0 0 5

# --- 3. Re-sort the country table (rows 4-216) by Casos totales desc -----
$sortRange = $ws.Range("A4:H216")
$keyRange = $ws.Range("B4:B216")
$sortRange.Sort($keyRange, 2, $null, $null, 1, $null, 1, 1, $false, $null, $null, 1)
